# Update cryptocurrency price/volume data to reflect latest scrape
# (values refreshed by the GitHub Actions job; two pairs of rows
# also swapped rank position: LidoDAOToken/Monero and Aptos/Frax)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.407.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.53%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.438.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.79%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.34%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "277.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3732"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3084"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.45"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.58%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.012"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06578"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -8.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.368"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.31"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -6.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.133"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -7.48%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.438.95"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -6.77%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001012"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -8.08%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "76.54"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -8.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05816"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -11.36%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.731"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.39"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.00%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.322"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.32%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.406.22"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.59%  "

# Row 26
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.231"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.02%  "

# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "142.48"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.03"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -7.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.602.08"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.64%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.33%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.907"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -19.62%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9127"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.470"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -7.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07708"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.368"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.85%  "

# Row 36
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.95"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.03%  "

# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.33%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05698"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.138"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.727"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.34%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1917"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.66%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02030"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.23%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.337"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -14.46%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.587"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5326"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.14"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5161"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.84%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.99"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.14%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.788"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.057"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.36%  "
